$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, matching style of existing header cells (bold, border, centered)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# time_taken values for each data row (2-134)
$timeValues = @(
    "2021-10-05 10:51:35.775037",
    "2021-10-05 10:51:35.775049",
    "2021-10-05 10:51:35.775052",
    "2021-10-05 10:51:35.775054",
    "2021-10-05 10:51:35.775057",
    "2021-10-05 10:51:35.775060",
    "2021-10-05 10:51:35.775063",
    "2021-10-05 10:51:35.775065",
    "2021-10-05 10:51:35.775068",
    "2021-10-05 10:51:35.775071",
    "2021-10-05 10:51:35.775073",
    "2021-10-05 10:51:35.775076",
    "2021-10-05 10:51:35.775078",
    "2021-10-05 10:51:35.775081",
    "2021-10-05 10:51:35.775083",
    "2021-10-05 10:51:35.775086",
    "2021-10-05 10:51:35.775088",
    "2021-10-05 10:51:35.775091",
    "2021-10-05 10:51:35.775094",
    "2021-10-05 10:51:35.775096",
    "2021-10-05 10:51:35.775099",
    "2021-10-05 10:51:35.775101",
    "2021-10-05 10:51:35.775104",
    "2021-10-05 10:51:35.775106",
    "2021-10-05 10:51:35.775109",
    "2021-10-05 10:51:35.775112",
    "2021-10-05 10:51:35.775114",
    "2021-10-05 10:51:35.775117",
    "2021-10-05 10:51:35.775119",
    "2021-10-05 10:51:35.775122",
    "2021-10-05 10:51:35.775124",
    "2021-10-05 10:51:35.775127",
    "2021-10-05 10:51:35.775130",
    "2021-10-05 10:51:35.775134",
    "2021-10-05 10:51:35.775137",
    "2021-10-05 10:51:35.775139",
    "2021-10-05 10:51:35.775142",
    "2021-10-05 10:51:35.775144",
    "2021-10-05 10:51:35.775147",
    "2021-10-05 10:51:35.775149",
    "2021-10-05 10:51:35.775153",
    "2021-10-05 10:51:35.775155",
    "2021-10-05 10:51:35.775158",
    "2021-10-05 10:51:35.775161",
    "2021-10-05 10:51:35.775163",
    "2021-10-05 10:51:35.775166",
    "2021-10-05 10:51:35.775168",
    "2021-10-05 10:51:35.775171",
    "2021-10-05 10:51:35.775173",
    "2021-10-05 10:51:35.775176",
    "2021-10-05 10:51:35.775178",
    "2021-10-05 10:51:35.775181",
    "2021-10-05 10:51:35.775184",
    "2021-10-05 10:51:35.775186",
    "2021-10-05 10:51:35.775189",
    "2021-10-05 10:51:35.775191",
    "2021-10-05 10:51:35.775194",
    "2021-10-05 10:51:35.775196",
    "2021-10-05 10:51:35.775199",
    "2021-10-05 10:51:35.775201",
    "2021-10-05 10:51:35.775204",
    "2021-10-05 10:51:35.775206",
    "2021-10-05 10:51:35.775209",
    "2021-10-05 10:51:35.775211",
    "2021-10-05 10:51:35.775216",
    "2021-10-05 10:51:35.775218",
    "2021-10-05 10:51:35.775221",
    "2021-10-05 10:51:35.775223",
    "2021-10-05 10:51:35.775226",
    "2021-10-05 10:51:35.775229",
    "2021-10-05 10:51:35.775231",
    "2021-10-05 10:51:35.775233",
    "2021-10-05 10:51:35.775236",
    "2021-10-05 10:51:35.775238",
    "2021-10-05 10:51:35.775241",
    "2021-10-05 10:51:35.775243",
    "2021-10-05 10:51:35.775248",
    "2021-10-05 10:51:35.775251",
    "2021-10-05 10:51:35.775254",
    "2021-10-05 10:51:35.775256",
    "2021-10-05 10:51:35.775259",
    "2021-10-05 10:51:35.775262",
    "2021-10-05 10:51:35.775264",
    "2021-10-05 10:51:35.775267",
    "2021-10-05 10:51:35.775269",
    "2021-10-05 10:51:35.775272",
    "2021-10-05 10:51:35.775274",
    "2021-10-05 10:51:35.775277",
    "2021-10-05 10:51:35.775279",
    "2021-10-05 10:51:35.775282",
    "2021-10-05 10:51:35.775285",
    "2021-10-05 10:51:35.775287",
    "2021-10-05 10:51:35.775291",
    "2021-10-05 10:51:35.775294",
    "2021-10-05 10:51:35.775296",
    "2021-10-05 10:51:35.775299",
    "2021-10-05 10:51:35.775301",
    "2021-10-05 10:51:35.775304",
    "2021-10-05 10:51:35.775306",
    "2021-10-05 10:51:35.775309",
    "2021-10-05 10:51:35.775311",
    "2021-10-05 10:51:35.775314",
    "2021-10-05 10:51:35.775316",
    "2021-10-05 10:51:35.775319",
    "2021-10-05 10:51:35.775321",
    "2021-10-05 10:51:35.775324",
    "2021-10-05 10:51:35.775326",
    "2021-10-05 10:51:35.775329",
    "2021-10-05 10:51:35.775334",
    "2021-10-05 10:51:35.775337",
    "2021-10-05 10:51:35.775339",
    "2021-10-05 10:51:35.775342",
    "2021-10-05 10:51:35.775345",
    "2021-10-05 10:51:35.775347",
    "2021-10-05 10:51:35.775350",
    "2021-10-05 10:51:35.775353",
    "2021-10-05 10:51:35.775355",
    "2021-10-05 10:51:35.775358",
    "2021-10-05 10:51:35.775360",
    "2021-10-05 10:51:35.775363",
    "2021-10-05 10:51:35.775366",
    "2021-10-05 10:51:35.775368",
    "2021-10-05 10:51:35.775371",
    "2021-10-05 10:51:35.775374",
    "2021-10-05 10:51:35.775376",
    "2021-10-05 10:51:35.775379",
    "2021-10-05 10:51:35.775381",
    "2021-10-05 10:51:35.775384",
    "2021-10-05 10:51:35.775388",
    "2021-10-05 10:51:35.775391",
    "2021-10-05 10:51:35.775394",
    "2021-10-05 10:51:35.775397",
    "2021-10-05 10:51:35.775399"
)

for ($i = 0; $i -lt $timeValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeValues[$i]
}

Write-Output "done"